# Update "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I) for rows 3-25
# on the "Training Dashboard" sheet: decrement H by 1 and set I to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 25; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $iCell = $ws.Cells.Item($row, 9)   # column I

    $hCell.Value2 = [double]$hCell.Value2 - 1
    # Leading apostrophe forces literal text so Excel doesn't auto-convert
    # the date-like string into a date serial number.
    $iCell.Value2 = "'04-Nov-2025"
}
